$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "67.747.62"
Set-TextValue "E2" "  +0.63%  "
Set-TextValue "D3" "3.814.44"
Set-TextValue "E3" "  +1.74%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "601.04"
Set-TextValue "E5" "  +1.29%  "
Set-TextValue "D6" "166.11"
Set-TextValue "E6" "  +0.37%  "
Set-TextValue "E7" "  -0.04%  "
Set-TextValue "E8" "  +0.23%  "
Set-TextValue "E9" "  +0.90%  "
Set-TextValue "D10" "6.31"
Set-TextValue "E10" "  -0.44%  "
Set-TextValue "E11" "  +0.88%  "
Set-TextValue "E12" "  -0.45%  "
Set-TextValue "D13" "35.74"
Set-TextValue "E13" "  -0.83%  "
Set-TextValue "D14" "4.457.50"
Set-TextValue "E14" "  +1.83%  "
Set-TextValue "D15" "3.815.56"
Set-TextValue "E15" "  +2.15%  "
Set-TextValue "D16" "67.781.46"
Set-TextValue "E16" "  +0.75%  "
Set-TextValue "D17" "18.40"
Set-TextValue "E17" "  +0.09%  "
Set-TextValue "D18" "7.06"
Set-TextValue "E18" "  +1.35%  "
Set-TextValue "E19" "  +0.52%  "
Set-TextValue "D20" "462.64"
Set-TextValue "E21" "  -0.84%  "
Set-TextValue "D22" "0.699"
Set-TextValue "E22" "  +0.70%  "
Set-TextValue "E23" "  -3.37%  "
Set-TextValue "D24" "83.34"
Set-TextValue "E24" "  +0.27%  "
Set-TextValue "E25" "  +1.78%  "
Set-TextValue "E26" "  -1.32%  "
Set-TextValue "D27" "10.05"
Set-TextValue "E27" "  -0.68%  "
Set-TextValue "E28" "  +0.01%  "
Set-TextValue "D29" "3.965.15"
Set-TextValue "E29" "  +1.76%  "
Set-TextValue "E30" "  +0.65%  "
Set-TextValue "D31" "7.37"
Set-TextValue "E31" "  +1.75%  "
Set-TextValue "E32" "  +2.53%  "
Set-TextValue "D33" "29.55"
Set-TextValue "E33" "  -0.05%  "
Set-TextValue "E34" "  -0.15%  "

Set-TextValue "B35" "Aptos"
Set-TextValue "C35" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D35" "9.08"
Set-TextValue "E35" "  -0.87%  "
Set-TextValue "B36" "Hedera"
Set-TextValue "C36" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D36" "0.0998"
Set-TextValue "E36" "  +0.03%  "
Set-TextValue "B37" "dogwifhat"
Set-TextValue "C37" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D37" "3.28"
Set-TextValue "E37" "  -0.92%  "
Set-TextValue "B38" "Kaspa"
Set-TextValue "C38" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D38" "0.138"
Set-TextValue "E38" "  +0.17%  "
Set-TextValue "B39" "Mantle"
Set-TextValue "C39" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D39" "0.996"
Set-TextValue "E39" "  +0.40%  "
Set-TextValue "B40" "Filecoin"
Set-TextValue "C40" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D40" "5.80"
Set-TextValue "E40" "  +1.11%  "
Set-TextValue "B41" "FirstDigitalUSD"
Set-TextValue "C41" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D41" "1.00"
Set-TextValue "E41" "  +0.11%  "
Set-TextValue "B42" "USDe"
Set-TextValue "C42" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "E42" "  -0.01%  "
Set-TextValue "B43" "OKB"
Set-TextValue "C43" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D43" "48.09"
Set-TextValue "E43" "  +2.36%  "
Set-TextValue "B44" "EnergySwap"
Set-TextValue "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "28.54"
Set-TextValue "E44" "  +9.38%  "
Set-TextValue "B45" "TheGraph"
Set-TextValue "C45" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D45" "0.300"
Set-TextValue "E45" "  +0.74%  "
Set-TextValue "B46" "Arweave"
Set-TextValue "C46" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D46" "43.22"
Set-TextValue "E46" "  -4.24%  "
Set-TextValue "B47" "ONDO"
Set-TextValue "C47" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D47" "1.40"
Set-TextValue "E47" "  +12.61%  "
Set-TextValue "B48" "Monero"
Set-TextValue "C48" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D48" "148.89"
Set-TextValue "E48" "  +0.24%  "
Set-TextValue "B49" "Cosmos"
Set-TextValue "C49" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D49" "8.33"
Set-TextValue "E49" "  +0.19%  "
Set-TextValue "B50" "Stacks"
Set-TextValue "C50" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D50" "1.83"
Set-TextValue "E50" "  +0.07%  "
Set-TextValue "B51" "Bittensor"
Set-TextValue "C51" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D51" "385.59"
Set-TextValue "E51" "  -0.97%  "
